# Add "SuperHide" to cell B1 on "Hoja3" and make that sheet the active one,
# with the selection moved to B2 (mirrors the author typing a value into
# B1 and then moving down to B2, leaving Hoja3 as the active/selected tab).

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Hoja3")
$ws3.Range("B1").Value = "SuperHide"

# Activating Hoja3 makes it the workbook's active tab and moves
# "tabSelected" away from whichever sheet had it before (Hoja4).
$ws3.Activate()
$ws3.Range("B2").Select()
